$wb = $excel.ActiveWorkbook

# --- Employees sheet: Jose Flores row (row 2) - hours C2: 0 -> "2" (text) ---
$wsEmployees = $wb.Worksheets.Item("Employees")
$wsEmployees.Range("C2").Value = "'2"
$wsEmployees.Range("C2").Style = "Normal"

# --- Equipment sheet: row 2 - swap equipment (Box Truck -> Paver) and update hours ---
$wsEquipment = $wb.Worksheets.Item("Equipment")
$wsEquipment.Range("A2").Value = "'712"
$wsEquipment.Range("A2").Style = "Normal"
$wsEquipment.Range("B2").Value = "Paver"
$wsEquipment.Range("C2").Value = "'2"
$wsEquipment.Range("C2").Style = "Normal"
